$d = $word.ActiveDocument

# 1. Merge the "rendez vous ... (réunion agile )" runs into a single run.
$d.Content.Find.Execute(
    "rendez vous à l’entreprise pour faire le point sur l’avancement des projets (réunion agile )",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "rendez vous à l’entreprise pour faire le point sur l’avancement des projets (réunion agile )",
    2) | Out-Null

# 2-4. Give the three empty paragraphs (between "soumission du cahier ... 1.0" and
# "Veille général sur les info de la tech et data" on "Jeudi 22 octobre") an explicit
# 11pt run size (sz/szCs = 22) matching their paragraph mark formatting.
$d.Paragraphs(220).Range.Font.Size = 11
$d.Paragraphs(220).Range.Font.SizeBi = 11
$d.Paragraphs(222).Range.Font.Size = 11
$d.Paragraphs(222).Range.Font.SizeBi = 11
$d.Paragraphs(224).Range.Font.Size = 11
$d.Paragraphs(224).Range.Font.SizeBi = 11

# 5. Remove the stray leading space run before "reprise leçon sur Microsoft learn".
$p226 = $d.Paragraphs(226)
$spaceRange = $d.Range($p226.Range.Start, $p226.Range.Start + 1)
$spaceRange.Delete()

# 6. Append the new "## après midi" afternoon section at the end of the document.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$afternoonPara = $d.Paragraphs($d.Paragraphs.Count)
$afternoonPara.Range.Text = "## après midi "

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$blankPara.Range.InsertParagraphAfter()

$lastBlank = $d.Paragraphs($d.Paragraphs.Count)
$lastBlank.Range.InsertParagraphAfter()
$finalPara = $d.Paragraphs($d.Paragraphs.Count)
$finalPara.Range.Text = "soumission des cahiers des charges avec toutes les phases réuni au Product Owner"
$finalPara.Range.Font.Size = 11
$finalPara.Range.Font.SizeBi = 11
